$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The review entry "Le nom de variable ""i"" n'est pas significatif." (old row 21,
# LIGNE "42") was dropped from the code-review table. Every row below it
# (old rows 22-33) shifts up by one (new rows 21-32), and each rows LIGNE
# (column B) value is updated to match the renumbered source file -- per
# "ligne changement selon document initiale". Columns C/D/E are unaffected
# by the shift, only LIGNE (B) changes for the shifted rows.

$ws.Cells.Item(12,2).Value = 7
$ws.Cells.Item(12,3).Value = "Erreur"
$ws.Cells.Item(12,4).Value = "Le chemin du fichier est mauvais."
$ws.Cells.Item(12,5).Value = "Correction"

$ws.Cells.Item(13,2).Value = "8 à 16"
$ws.Cells.Item(13,3).Value = "Amélioration"
$ws.Cells.Item(13,4).Value = "Mettre ce bloc de code dans une fonction."
$ws.Cells.Item(13,5).Value = "Billet"

$ws.Cells.Item(14,2).Value = "7 à 110"
$ws.Cells.Item(14,3).Value = "Amélioration"
$ws.Cells.Item(14,4).Value = "Mettre dans le main (incluant les fonctions)."
$ws.Cells.Item(14,5).Value = "Billet"

$ws.Cells.Item(15,2).Value = 14
$ws.Cells.Item(15,3).Value = "Cosmétique"
$ws.Cells.Item(15,4).Value = "La variable fichier est en français alors que la majorité des variables sont en anglais."
$ws.Cells.Item(15,5).Value = "Correction"

$ws.Cells.Item(16,2).Value = 23
$ws.Cells.Item(16,3).Value = "Erreur"
$ws.Cells.Item(16,4).Value = "La variable choice est initialisé à 0, pourquoi, si c'est un input de caractères. Initialiser à vide."
$ws.Cells.Item(16,5).Value = "Correction"

$ws.Cells.Item(17,2).Value = "33 à 41"
$ws.Cells.Item(17,3).Value = "Amélioration"
$ws.Cells.Item(17,4).Value = "Mettre ce bloc de code dans une fonction."
$ws.Cells.Item(17,5).Value = "Billet"

$ws.Cells.Item(18,2).Value = 33
$ws.Cells.Item(18,3).Value = "Cosmétique"
$ws.Cells.Item(18,4).Value = "Le nom de variable ""indexes"" n'est pas significatif."
$ws.Cells.Item(18,5).Value = "Correction"

$ws.Cells.Item(19,2).Value = 34
$ws.Cells.Item(19,3).Value = "Amélioration"
$ws.Cells.Item(19,4).Value = "La variable ""i"" dans le for n'est pas utilisé. Mettre à la place ""_""."
$ws.Cells.Item(19,5).Value = "Correction"

$ws.Cells.Item(20,2).Value = 35
$ws.Cells.Item(20,3).Value = "Cosmétique"
$ws.Cells.Item(20,4).Value = "Le nom de variable""index"" n'est pas significatif. "
$ws.Cells.Item(20,5).Value = "Correction"

$ws.Cells.Item(21,2).Value = 45
$ws.Cells.Item(21,3).Value = "Cosmétique"
$ws.Cells.Item(21,4).Value = "La variable fils_coupes est en français alors que la majorité des variables sont en anglais."
$ws.Cells.Item(21,5).Value = "Correction"

$ws.Cells.Item(22,2).Value = 47
$ws.Cells.Item(22,3).Value = "Erreur"
$ws.Cells.Item(22,4).Value = "La variable ""fin"" devrait être un boolean."
$ws.Cells.Item(22,5).Value = "Correction"

$ws.Cells.Item(23,2).Value = 47
$ws.Cells.Item(23,3).Value = "Cosmétique"
$ws.Cells.Item(23,4).Value = "La variable fin est en français alors que la majorité des variables sont en anglais."
$ws.Cells.Item(23,5).Value = "Correction"

$ws.Cells.Item(24,2).Value = 48
$ws.Cells.Item(24,3).Value = "Cosmétique"
$ws.Cells.Item(24,4).Value = "Le nom de variable ""n"" est non-significatif."
$ws.Cells.Item(24,5).Value = "Correction"

$ws.Cells.Item(25,2).Value = "De 50 à 101"
$ws.Cells.Item(25,3).Value = "Amélioration"
$ws.Cells.Item(25,4).Value = "Manque de commentaires qui assurent la compréhension du code."
$ws.Cells.Item(25,5).Value = "Billet"

$ws.Cells.Item(26,2).Value = "51 à 55"
$ws.Cells.Item(26,3).Value = "Amélioration"
$ws.Cells.Item(26,4).Value = "Mettre ce bloc de code dans une fonction."
$ws.Cells.Item(26,5).Value = "Correction"

$ws.Cells.Item(27,2).Value = 57
$ws.Cells.Item(27,3).Value = "Amélioration"
$ws.Cells.Item(27,4).Value = "Il serait mieux de faire un retour de ligne au lieu d'un print() vide."
$ws.Cells.Item(27,5).Value = "Correction"

$ws.Cells.Item(28,2).Value = 59
$ws.Cells.Item(28,3).Value = "Cosmétique"
$ws.Cells.Item(28,4).Value = "La variable fil est en français alors que la majorité des variables sont en anglais."
$ws.Cells.Item(28,5).Value = "Correction"

$ws.Cells.Item(29,2).Value = "61 à 101 "
$ws.Cells.Item(29,3).Value = "Amélioration"
$ws.Cells.Item(29,4).Value = "Mettre ce bloc de code dans une fonction."
$ws.Cells.Item(29,5).Value = "Billet"

$ws.Cells.Item(30,2).Value = 59
$ws.Cells.Item(30,3).Value = "Erreur"
$ws.Cells.Item(30,4).Value = "L'input ""fil"" devrait obliger le joueur à entrer une valeur entre 1 et 5."
$ws.Cells.Item(30,5).Value = "Billet"

$ws.Cells.Item(31,2).Value = 90
$ws.Cells.Item(31,3).Value = "Erreur"
$ws.Cells.Item(31,4).Value = "Le traitement dans le cas de fil = 5 est manquant."
$ws.Cells.Item(31,5).Value = "Billet"

$ws.Cells.Item(32,2).Value = "103 à 106"
$ws.Cells.Item(32,3).Value = "Amélioration"
$ws.Cells.Item(32,4).Value = "Mettre ce bloc de code dans une fonction."
$ws.Cells.Item(32,5).Value = "Billet"

# Drop the now-superfluous last row (data now ends at row 32); clearing the
# full row removes it from the saved sheetData and shrinks the used range.
$ws.Range("A33:G33").Clear()

# Match the saved view: scrolled down one row, with B32 as the active selection
$ws.Range("B32").Select()
